$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdFile = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
$zhXlf = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
$deXlf = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/faeab19840ee0b595bfee01571bffad4b7f7cfd7/e2e/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/05d194503093fbd283883c17b37f0ce33af4a026/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/070c20bb87cfd1cb0a7a19263edcd05aa8609905/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"

# 1. Update the "Ready for handoff" status everywhere it appears so the
#    shared string itself now reads "Handed back: in sync with en-US".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

# 2. Fill in the "Latest Target File" (E) / "Latest Handback File" (F)
#    columns for the data rows, and stamp the real handback datetime into
#    "Latest Handback DateTime" (G), now that the handback report exists.
$wsZh.Range("E2").Value = $mdFile
$wsZh.Range("F2").Value = $zhXlf
$wsZh.Range("G2").Value = "2016-03-10 01:16:20"

$wsZh.Range("E3").Value = $mdFile
$wsZh.Range("F3").Value = $zhXlf
$wsZh.Range("G3").Value = "2016-03-10 01:16:20"

$wsDe.Range("E2").Value = $mdFile
$wsDe.Range("F2").Value = $deXlf
$wsDe.Range("G2").Value = "2016-03-10 01:16:38"

$wsDe.Range("E3").Value = $mdFile
$wsDe.Range("F3").Value = $deXlf
$wsDe.Range("G3").Value = "2016-03-10 01:16:38"

# 3. Hyperlink the newly populated E/F cells to the same targets as their
#    A/C counterparts. Existing hyperlinks (A2, C2, A3, C3, A4) are left
#    untouched so their styling doesn't get disturbed.
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl, "", "", $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlfUrl, "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $mdUrl, "", "", $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhXlfUrl, "", "", $zhXlf)

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl, "", "", $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlfUrl, "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $mdUrl, "", "", $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deXlfUrl, "", "", $deXlf)

Write-Host "Generated handback report"
